$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price column D, Volume(1h) column E).
# Values are written as text (matching the source data format) by
# prefixing with a leading apostrophe, then the cell style is reset
# to Normal so no stray "quote prefix" formatting is left behind.

$c = $ws.Range("D2")
$c.Value = "'334.82"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "'1.62%"
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'43.98"
$c.Style = "Normal"

$c = $ws.Range("E3")
$c.Value = "'6.63%"
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'5.742"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'1.70%"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'0.08338"
$c.Style = "Normal"

$c = $ws.Range("E5")
$c.Value = "'1.54%"
$c.Style = "Normal"

$c = $ws.Range("D6")
$c.Value = "'8.845"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'1.00%"
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'4.519"
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'0.36%"
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'1.963"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'-2.65%"
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'2.879"
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "'-3.60%"
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.9506"
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'3.05%"
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.1249"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'-2.16%"
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'0.1985"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'1.31%"
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'0.1056"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'12.78%"
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'0.04536"
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'18.00%"
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.1067"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'0.79%"
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'0.001296"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'-0.76%"
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'0.005970"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'-2.52%"
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'3.499"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'1.49%"
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'8.688"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'5.05%"
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'0.1362"
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'-0.28%"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'1.16%"
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'0.04413"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'0.20%"
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'0.001257"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'-0.12%"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'0.004336"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'0.50%"
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'0.0001262"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'5.11%"
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'0.0003993"
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.02810"
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "'1.91%"
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.06038"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'10.44%"
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'0.007919"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'-1.18%"
$c.Style = "Normal"

$c = $ws.Range("E42")
$c.Value = "'0.52%"
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'0.008975"
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "'0.36%"
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'0.002145"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'-1.19%"
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'-11.45%"
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'0.00007008"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'0.10%"
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'0.003187"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'-0.09%"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'-0.35%"
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'0.10%"
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'0.0002003"
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'0.10%"
$c.Style = "Normal"

